# Daily attendance processing - 2026-01-24 13:47:12
# Normalize the "Recorded By" column (G) so that the author's email
# address is listed before "System" instead of after it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
